# PlayerPerformance_4825.xlsx update
# - Adds a new "Player Info" sheet (before "ODI Batting") with the
#   player's basic info (ID / NAME / BATTING_HAND / BOWL_STYLE).
# - Renames MATCH_CARD_LINK -> MATCH_CODE on "ODI Batting" (col D) and
#   "ODI Bowling" (col B), and replaces the full howstat.com scorecard
#   URL with the bare numeric match code it encoded.

$wb = $excel.ActiveWorkbook

# --- 1. Insert a new "Player Info" worksheet right before "ODI Batting" ---
# NOTE: worksheet references returned by this host resolve by *position*,
# not by a stable object handle, so once a new sheet is spliced in front
# of "ODI Batting" any variable captured beforehand now points at the
# new sheet instead. Grab a fresh reference by name for the insertion
# point, then do the same for "ODI Batting"/"ODI Bowling" again after
# the insertion has shifted everybody's index.
$insertBeforeMe = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($insertBeforeMe)
$playerInfo.Name = "Player Info"

# Header row
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Match the bold / centered / bordered header look used on the other sheets
$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108  # xlCenter
$headerRange.VerticalAlignment = -4160    # xlTop
$headerRange.Borders.LineStyle = 1        # xlContinuous

# Data row. Keep the ID as text (matches how the match codes below, and
# every other "numberish" column in this workbook, are stored).
$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "4825"
$playerInfo.Range("B2").Value = "Obed Christopher McCoy"
$playerInfo.Range("C2").Value = "Left Handed"
$playerInfo.Range("D2").Value = "Left Arm Fast Medium"

# --- 2. "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE (column D) ---
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("D1").Value = "MATCH_CODE"

$battingSheet.Range("D2").NumberFormat = "@"
$battingSheet.Range("D2").Value = "4216"

$battingSheet.Range("D3").NumberFormat = "@"
$battingSheet.Range("D3").Value = "4219"

# --- 3. "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE (column B) ---
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

$bowlingSheet.Range("B2").NumberFormat = "@"
$bowlingSheet.Range("B2").Value = "4216"

$bowlingSheet.Range("B3").NumberFormat = "@"
$bowlingSheet.Range("B3").Value = "4219"

Write-Output "Sheets: $([string]::Join(', ', ($wb.Worksheets | ForEach-Object { $_.Name })))"
